# Estadisticos Segundo Parcial 23 Mayo
$wb = $excel.ActiveWorkbook

# --- Estadisticos 2P (sheet2): update Blancos/Reprobados/Aprobados/Por_Apro/Promedio ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$data2P = @(
    @(2, 36, 0, 14, 22, 61.11, 5.7),
    @(3, 28, 0, 6, 22, 78.57, 6.2),
    @(4, 30, 1, 5, 25, 83.33, 7)
)
foreach ($row in $data2P) {
    $r = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $ws2.Cells.Item($r, 6).Value = $row[4]
    $ws2.Cells.Item($r, 7).Value = $row[5]
    $ws2.Cells.Item($r, 8).Value = $row[6]
}

# --- Estadisticos Final (sheet3): update Reprobados/Aprobados/Por_Apro/Promedio ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$data3P = @(
    @(2, 0, 14, 22, 61.11, 6.6),
    @(3, 0, 6, 22, 78.57, 7.3),
    @(4, 0, 5, 25, 83.33, 7.4)
)
foreach ($row in $data3P) {
    $r = $row[0]
    $ws3.Cells.Item($r, 4).Value = $row[1]
    $ws3.Cells.Item($r, 5).Value = $row[2]
    $ws3.Cells.Item($r, 6).Value = $row[3]
    $ws3.Cells.Item($r, 7).Value = $row[4]
    $ws3.Cells.Item($r, 8).Value = $row[5]
}

# --- Rescatables (sheet4): refresh roster of at-risk students (rows 2-20) ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$dataResc = @(
    @(2, 24330051920092, "APARICIO", "OFICIAL", "VICTOR YAEL", "Pensamiento matemático II", "2AEV", 4),
    @(3, 24330051920304, "ARMAS", "SALINAS", "JOSE GUSTAVO", "Pensamiento matemático II", "2AEV", 4),
    @(4, 24330051920305, "MORALES", "CUAHUA", "ANDRES", "Pensamiento matemático II", "2AEV", 4),
    @(5, 24330051920113, "RAMOS", "DE LA CRUZ", "DEREK", "Pensamiento matemático II", "2AEV", 4),
    @(6, 24330051920315, "VENTURA", "ZEPEDA", "CARLOS ARGEL", "Pensamiento matemático II", "2AEV", 4),
    @(7, 24330051920389, "RUIZ", "MORALES", "MAYRIN GUADALUPE", "Pensamiento matemático II", "2ALCV", 4),
    @(8, 24330051920330, "VASQUEZ", "PEREZ", "DANIELA LILI", "Pensamiento matemático II", "2ALCV", 4),
    @(9, 24330051920246, "ZUNO", "FLORES", "ALIN MARIEL", "Pensamiento matemático II", "2ALCV", 4),
    @(10, 24330051920182, "LOPEZ", "DE LA CRUZ", "AMISADAY", "Pensamiento matemático II", "2ARHV", 4),
    @(11, 24330051920093, "ARIAS", "SARMIENTO", "URIEL ARTURO", "Pensamiento matemático II", "2AEV", 3),
    @(12, 24330051920090, "ANTONIO", "LOPEZ", "SERGIO GISELL", "Pensamiento matemático II", "2AEV", 3),
    @(13, 24330051920098, "CHICO", "BALDERAS", "YARETH", "Pensamiento matemático II", "2AEV", 3),
    @(14, 23330051920036, "HERNANDEZ", "DOLORES", "GONZALO", "Pensamiento matemático II", "2AEV", 3),
    @(15, 24330051920144, "MUÑOZ", "CORONA", "JOSE ABEL", "Pensamiento matemático II", "2AEV", 3),
    @(16, 24330051920306, "ROJAS", "GUTIERREZ", "LUIS ROBERTO", "Pensamiento matemático II", "2AEV", 3),
    @(17, 24330051920238, "TORRES", "PEREZ", "ERIKA VALERIA", "Pensamiento matemático II", "2ALCV", 3),
    @(18, 23330051920298, "MAZA", "ENCARNACION", "KEVIN JESUS", "Pensamiento matemático II", "2ARHV", 3),
    @(19, 23330051920224, "DORANTES", "PORRAS", "ROBERTO", "Pensamiento matemático II", "2AEV", 2),
    @(20, 24330051920206, "PORTUGAL", "VEGA", "SANTIAGO", "Pensamiento matemático II", "2ARHV", 2)
)
foreach ($row in $dataResc) {
    $r = $row[0]
    $ws4.Cells.Item($r, 1).Value = $row[1]
    $ws4.Cells.Item($r, 2).Value = $row[2]
    $ws4.Cells.Item($r, 3).Value = $row[3]
    $ws4.Cells.Item($r, 4).Value = $row[4]
    $ws4.Cells.Item($r, 5).Value = $row[5]
    $ws4.Cells.Item($r, 6).Value = $row[6]
    $ws4.Cells.Item($r, 7).Value = $row[7]
}

Write-Host "Estadisticos actualizados: 2P, Final y Rescatables"
